$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("H:H").Insert()
$ws.Range("H1").Value = "Hire Date"
$ws.Range("H2").Select()
